# Column annotation unit testcase implemented
#
# 1. Remove the unused blank "Sheet1" tab.
# 2. Leave cell A4 selected on the "User" sheet (its last recorded selection).
# 3. Add a new "UserWithAnnotations" sheet at the end of the workbook that
#    mirrors the "User" sheet's data, but with its second column header
#    renamed to "fullName" (exercising a column-annotation/mapping test).
# 4. Make the new sheet the active tab.

$wb = $excel.ActiveWorkbook

# --- 1. Delete the empty "Sheet1" tab -------------------------------------
$wb.Worksheets.Item("Sheet1").Delete()

# --- 2. Record the A4 selection on the "User" sheet -----------------------
$userSheet = $wb.Worksheets.Item("User")
$userSheet.Range("A4").Select()

# --- 3. Add the new "UserWithAnnotations" sheet at the end ----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "UserWithAnnotations"

$newSheet.Range("A1").Value = "id"
$newSheet.Range("B1").Value = "fullName"

$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "One"

$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "Two"

$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "Three"

# --- 4. Activate the new sheet so it becomes the selected tab -------------
$newSheet.Activate()
$newSheet.Range("A1").Select()
